$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90
$ws.Cells.Item(90, 2).Value = 6978747
$ws.Cells.Item(90, 5).Value = "IMT Novi Belgrade"
$ws.Cells.Item(90, 6).Value = "Red Star Belgrade"
$ws.Cells.Item(90, 8).Value = 2
$ws.Cells.Item(90, 10).Value = 2
$ws.Cells.Item(90, 11).Value = "A"
$ws.Cells.Item(90, 12).Value = 8
$ws.Cells.Item(90, 13).Value = 5.25
$ws.Cells.Item(90, 14).Value = 1.285
$ws.Cells.Item(90, 15).Value = 15
$ws.Cells.Item(90, 16).Value = 7.5
$ws.Cells.Item(90, 17).Value = 1.125
$ws.Cells.Item(90, 18).Value = 2.25
$ws.Cells.Item(90, 19).Value = 1.975
$ws.Cells.Item(90, 20).Value = 1.825
$ws.Cells.Item(90, 21).Value = 3.5
$ws.Cells.Item(90, 22).Value = 1.825
$ws.Cells.Item(90, 23).Value = 1.975
$ws.Cells.Item(90, 25).Value = -1
$ws.Cells.Item(90, 26).Value = 0.125
$ws.Cells.Item(90, 27).Value = 0.9750000000000001
$ws.Cells.Item(90, 28).Value = -1
$ws.Cells.Item(90, 30).Value = 0.9750000000000001

# Row 91
$ws.Cells.Item(91, 2).Value = 6979491
$ws.Cells.Item(91, 5).Value = "Radnicki Nis"
$ws.Cells.Item(91, 6).Value = "Spartak Subotica"
$ws.Cells.Item(91, 8).Value = 1
$ws.Cells.Item(91, 10).Value = 1
$ws.Cells.Item(91, 11).Value = "D"
$ws.Cells.Item(91, 12).Value = 1.95
$ws.Cells.Item(91, 13).Value = 3.25
$ws.Cells.Item(91, 14).Value = 3.7
$ws.Cells.Item(91, 15).Value = 1.65
$ws.Cells.Item(91, 16).Value = 3.5
$ws.Cells.Item(91, 17).Value = 5
$ws.Cells.Item(91, 18).Value = -0.75
$ws.Cells.Item(91, 19).Value = 1.825
$ws.Cells.Item(91, 20).Value = 1.975
$ws.Cells.Item(91, 21).Value = 2.5
$ws.Cells.Item(91, 22).Value = 2
$ws.Cells.Item(91, 23).Value = 1.8
$ws.Cells.Item(91, 25).Value = 2.5
$ws.Cells.Item(91, 26).Value = -1
$ws.Cells.Item(91, 27).Value = -1
$ws.Cells.Item(91, 28).Value = 0.9750000000000001
$ws.Cells.Item(91, 30).Value = 0.8

# Row 239
$ws.Cells.Item(239, 2).Value = 6979610
$ws.Cells.Item(239, 5).Value = "Spartak Subotica"
$ws.Cells.Item(239, 6).Value = "Javor Ivanjica"
$ws.Cells.Item(239, 8).Value = 0
$ws.Cells.Item(239, 9).Value = 0
$ws.Cells.Item(239, 11).Value = "H"
$ws.Cells.Item(239, 12).Value = 2.1
$ws.Cells.Item(239, 13).Value = 3.2
$ws.Cells.Item(239, 14).Value = 3.3
$ws.Cells.Item(239, 15).Value = 2
$ws.Cells.Item(239, 16).Value = 3.3
$ws.Cells.Item(239, 17).Value = 3.5
$ws.Cells.Item(239, 18).Value = -0.5
$ws.Cells.Item(239, 21).Value = 2.25
$ws.Cells.Item(239, 22).Value = 1.775
$ws.Cells.Item(239, 23).Value = 2.025
$ws.Cells.Item(239, 24).Value = 1
$ws.Cells.Item(239, 26).Value = -1
$ws.Cells.Item(239, 27).Value = 1
$ws.Cells.Item(239, 28).Value = -1
$ws.Cells.Item(239, 29).Value = -1
$ws.Cells.Item(239, 30).Value = 1.025

# Row 240
$ws.Cells.Item(240, 2).Value = 6979611
$ws.Cells.Item(240, 5).Value = "Mladost Lucani"
$ws.Cells.Item(240, 6).Value = "IMT Novi Belgrade"
$ws.Cells.Item(240, 8).Value = 2
$ws.Cells.Item(240, 9).Value = 1
$ws.Cells.Item(240, 11).Value = "A"
$ws.Cells.Item(240, 12).Value = 2.3
$ws.Cells.Item(240, 13).Value = 3.4
$ws.Cells.Item(240, 14).Value = 2.75
$ws.Cells.Item(240, 15).Value = 2.9
$ws.Cells.Item(240, 16).Value = 4
$ws.Cells.Item(240, 17).Value = 2
$ws.Cells.Item(240, 18).Value = 0.25
$ws.Cells.Item(240, 21).Value = 3
$ws.Cells.Item(240, 22).Value = 2
$ws.Cells.Item(240, 23).Value = 1.8
$ws.Cells.Item(240, 24).Value = -1
$ws.Cells.Item(240, 26).Value = 1
$ws.Cells.Item(240, 27).Value = -1
$ws.Cells.Item(240, 28).Value = 0.8
$ws.Cells.Item(240, 29).Value = 0
$ws.Cells.Item(240, 30).Value = 0

# Row 267
$ws.Cells.Item(267, 2).Value = 8106961
$ws.Cells.Item(267, 5).Value = "Spartak Subotica"
$ws.Cells.Item(267, 6).Value = "Radnicki Nis"
$ws.Cells.Item(267, 7).Value = 3
$ws.Cells.Item(267, 9).Value = 3
$ws.Cells.Item(267, 11).Value = "H"
$ws.Cells.Item(267, 12).Value = 2.4
$ws.Cells.Item(267, 14).Value = 2.625
$ws.Cells.Item(267, 15).Value = 3
$ws.Cells.Item(267, 16).Value = 2.15
$ws.Cells.Item(267, 17).Value = 3.3
$ws.Cells.Item(267, 18).Value = 0
$ws.Cells.Item(267, 19).Value = 1.8
$ws.Cells.Item(267, 20).Value = 2
$ws.Cells.Item(267, 21).Value = 2
$ws.Cells.Item(267, 22).Value = 1.85
$ws.Cells.Item(267, 23).Value = 1.95
$ws.Cells.Item(267, 24).Value = 2
$ws.Cells.Item(267, 26).Value = -1
$ws.Cells.Item(267, 27).Value = 0.8
$ws.Cells.Item(267, 28).Value = -1
$ws.Cells.Item(267, 29).Value = 0.8500000000000001
$ws.Cells.Item(267, 30).Value = -1

# Row 268
$ws.Cells.Item(268, 2).Value = 8105038
$ws.Cells.Item(268, 5).Value = "FK Radnik Surdulica"
$ws.Cells.Item(268, 6).Value = "Javor Ivanjica"
$ws.Cells.Item(268, 7).Value = 0
$ws.Cells.Item(268, 9).Value = 0
$ws.Cells.Item(268, 11).Value = "A"
$ws.Cells.Item(268, 12).Value = 2.5
$ws.Cells.Item(268, 14).Value = 2.5
$ws.Cells.Item(268, 15).Value = 3.5
$ws.Cells.Item(268, 16).Value = 3.3
$ws.Cells.Item(268, 17).Value = 1.85
$ws.Cells.Item(268, 18).Value = 0.5
$ws.Cells.Item(268, 19).Value = 1.875
$ws.Cells.Item(268, 20).Value = 1.925
$ws.Cells.Item(268, 21).Value = 2.25
$ws.Cells.Item(268, 22).Value = 1.775
$ws.Cells.Item(268, 23).Value = 2.025
$ws.Cells.Item(268, 24).Value = -1
$ws.Cells.Item(268, 26).Value = 0.8500000000000001
$ws.Cells.Item(268, 27).Value = -1
$ws.Cells.Item(268, 28).Value = 0.925
$ws.Cells.Item(268, 29).Value = -1
$ws.Cells.Item(268, 30).Value = 1.025

# Row 279
$ws.Cells.Item(279, 2).Value = 8106964
$ws.Cells.Item(279, 5).Value = "Radnicki Nis"
$ws.Cells.Item(279, 6).Value = "IMT Novi Belgrade"
$ws.Cells.Item(279, 7).Value = 0
$ws.Cells.Item(279, 8).Value = 0
$ws.Cells.Item(279, 9).Value = 0
$ws.Cells.Item(279, 10).Value = 0
$ws.Cells.Item(279, 12).Value = 2.1
$ws.Cells.Item(279, 13).Value = 3.1
$ws.Cells.Item(279, 14).Value = 3.2
$ws.Cells.Item(279, 16).Value = 3.3
$ws.Cells.Item(279, 17).Value = 2.9
$ws.Cells.Item(279, 19).Value = 1.925
$ws.Cells.Item(279, 20).Value = 1.875
$ws.Cells.Item(279, 21).Value = 2.5
$ws.Cells.Item(279, 25).Value = 2.3
$ws.Cells.Item(279, 28).Value = 0.4375
$ws.Cells.Item(279, 29).Value = -1
$ws.Cells.Item(279, 30).Value = 0.8999999999999999

# Row 280
$ws.Cells.Item(280, 2).Value = 8106767
$ws.Cells.Item(280, 5).Value = "Javor Ivanjica"
$ws.Cells.Item(280, 6).Value = "FK Vozdovac"
$ws.Cells.Item(280, 7).Value = 1
$ws.Cells.Item(280, 8).Value = 1
$ws.Cells.Item(280, 9).Value = 1
$ws.Cells.Item(280, 10).Value = 1
$ws.Cells.Item(280, 12).Value = 2.25
$ws.Cells.Item(280, 13).Value = 3
$ws.Cells.Item(280, 14).Value = 3
$ws.Cells.Item(280, 16).Value = 3
$ws.Cells.Item(280, 17).Value = 3.1
$ws.Cells.Item(280, 19).Value = 1.9
$ws.Cells.Item(280, 20).Value = 1.9
$ws.Cells.Item(280, 21).Value = 2.25
$ws.Cells.Item(280, 25).Value = 2
$ws.Cells.Item(280, 28).Value = 0.45
$ws.Cells.Item(280, 29).Value = -0.5
$ws.Cells.Item(280, 30).Value = 0.45

# Row 284
$ws.Cells.Item(284, 2).Value = 8105865
$ws.Cells.Item(284, 5).Value = "Mladost Lucani"
$ws.Cells.Item(284, 6).Value = "FK Backa Topola"
$ws.Cells.Item(284, 7).Value = 2
$ws.Cells.Item(284, 8).Value = 0
$ws.Cells.Item(284, 9).Value = 1
$ws.Cells.Item(284, 11).Value = "H"
$ws.Cells.Item(284, 12).Value = 3.75
$ws.Cells.Item(284, 13).Value = 3.75
$ws.Cells.Item(284, 14).Value = 1.727
$ws.Cells.Item(284, 15).Value = 5.25
$ws.Cells.Item(284, 16).Value = 4.5
$ws.Cells.Item(284, 17).Value = 1.42
$ws.Cells.Item(284, 18).Value = 1.25
$ws.Cells.Item(284, 19).Value = 1.9
$ws.Cells.Item(284, 20).Value = 1.9
$ws.Cells.Item(284, 22).Value = 1.9
$ws.Cells.Item(284, 23).Value = 1.9
$ws.Cells.Item(284, 24).Value = 4.25
$ws.Cells.Item(284, 26).Value = -1
$ws.Cells.Item(284, 27).Value = 0.8999999999999999
$ws.Cells.Item(284, 28).Value = -1
$ws.Cells.Item(284, 30).Value = 0.8999999999999999

# Row 285
$ws.Cells.Item(285, 2).Value = 8105026
$ws.Cells.Item(285, 5).Value = "FK Cukaricki"
$ws.Cells.Item(285, 6).Value = "Partizan Belgrade"
$ws.Cells.Item(285, 7).Value = 0
$ws.Cells.Item(285, 8).Value = 1
$ws.Cells.Item(285, 9).Value = 0
$ws.Cells.Item(285, 11).Value = "A"
$ws.Cells.Item(285, 12).Value = 2.75
$ws.Cells.Item(285, 13).Value = 3.2
$ws.Cells.Item(285, 14).Value = 2.3
$ws.Cells.Item(285, 15).Value = 2.45
$ws.Cells.Item(285, 16).Value = 3.6
$ws.Cells.Item(285, 17).Value = 2.45
$ws.Cells.Item(285, 18).Value = 0
$ws.Cells.Item(285, 19).Value = 1.875
$ws.Cells.Item(285, 20).Value = 1.925
$ws.Cells.Item(285, 22).Value = 1.8
$ws.Cells.Item(285, 23).Value = 2
$ws.Cells.Item(285, 24).Value = -1
$ws.Cells.Item(285, 26).Value = 1.45
$ws.Cells.Item(285, 27).Value = -1
$ws.Cells.Item(285, 28).Value = 0.925
$ws.Cells.Item(285, 30).Value = 1

# Row 291
$ws.Cells.Item(291, 2).Value = 8245724
$ws.Cells.Item(291, 5).Value = "FK Novi Pazar"
$ws.Cells.Item(291, 6).Value = "Radnicki Nis"
$ws.Cells.Item(291, 7).Value = 2
$ws.Cells.Item(291, 8).Value = 1
$ws.Cells.Item(291, 9).Value = 1
$ws.Cells.Item(291, 12).Value = 4
$ws.Cells.Item(291, 13).Value = 2.6
$ws.Cells.Item(291, 14).Value = 2.1
$ws.Cells.Item(291, 15).Value = 2.625
$ws.Cells.Item(291, 16).Value = 2.8
$ws.Cells.Item(291, 17).Value = 2.8
$ws.Cells.Item(291, 18).Value = 0
$ws.Cells.Item(291, 19).Value = 1.825
$ws.Cells.Item(291, 20).Value = 1.975
$ws.Cells.Item(291, 21).Value = 2.5
$ws.Cells.Item(291, 24).Value = 1.625
$ws.Cells.Item(291, 27).Value = 0.825
$ws.Cells.Item(291, 28).Value = -1

# Row 292
$ws.Cells.Item(292, 2).Value = 8245726
$ws.Cells.Item(292, 5).Value = "Spartak Subotica"
$ws.Cells.Item(292, 6).Value = "IMT Novi Belgrade"
$ws.Cells.Item(292, 12).Value = 4.4
$ws.Cells.Item(292, 13).Value = 2.75
$ws.Cells.Item(292, 14).Value = 1.909
$ws.Cells.Item(292, 15).Value = 4
$ws.Cells.Item(292, 16).Value = 2.875
$ws.Cells.Item(292, 17).Value = 1.95
$ws.Cells.Item(292, 18).Value = 0.5
$ws.Cells.Item(292, 19).Value = 1.8
$ws.Cells.Item(292, 20).Value = 2
$ws.Cells.Item(292, 22).Value = 1.95
$ws.Cells.Item(292, 23).Value = 1.85
$ws.Cells.Item(292, 25).Value = 1.875
$ws.Cells.Item(292, 27).Value = 0.8
$ws.Cells.Item(292, 28).Value = -1
$ws.Cells.Item(292, 30).Value = 0.425

# Row 293
$ws.Cells.Item(293, 2).Value = 8245727
$ws.Cells.Item(293, 5).Value = "FK Vozdovac"
$ws.Cells.Item(293, 6).Value = "FK Radnik Surdulica"
$ws.Cells.Item(293, 7).Value = 3
$ws.Cells.Item(293, 8).Value = 2
$ws.Cells.Item(293, 10).Value = 1
$ws.Cells.Item(293, 11).Value = "H"
$ws.Cells.Item(293, 12).Value = 1.2
$ws.Cells.Item(293, 13).Value = 6.5
$ws.Cells.Item(293, 14).Value = 8
$ws.Cells.Item(293, 15).Value = 1.38
$ws.Cells.Item(293, 16).Value = 5.5
$ws.Cells.Item(293, 17).Value = 5.25
$ws.Cells.Item(293, 18).Value = -1.25
$ws.Cells.Item(293, 19).Value = 1.85
$ws.Cells.Item(293, 20).Value = 1.95
$ws.Cells.Item(293, 21).Value = 2.75
$ws.Cells.Item(293, 22).Value = 1.85
$ws.Cells.Item(293, 23).Value = 1.95
$ws.Cells.Item(293, 24).Value = 0.3799999999999999
$ws.Cells.Item(293, 25).Value = -1
$ws.Cells.Item(293, 27).Value = -0.5
$ws.Cells.Item(293, 28).Value = 0.475
$ws.Cells.Item(293, 29).Value = 0.8500000000000001
$ws.Cells.Item(293, 30).Value = -1

# Row 294
$ws.Cells.Item(294, 2).Value = 8245725
$ws.Cells.Item(294, 5).Value = "Javor Ivanjica"
$ws.Cells.Item(294, 6).Value = "FK Zeleznicar Pancevo"
$ws.Cells.Item(294, 7).Value = 1
$ws.Cells.Item(294, 9).Value = 0
$ws.Cells.Item(294, 10).Value = 0
$ws.Cells.Item(294, 11).Value = "D"
$ws.Cells.Item(294, 12).Value = 2.375
$ws.Cells.Item(294, 13).Value = 2.875
$ws.Cells.Item(294, 14).Value = 2.9
$ws.Cells.Item(294, 15).Value = 3.8
$ws.Cells.Item(294, 16).Value = 2.15
$ws.Cells.Item(294, 17).Value = 2.625
$ws.Cells.Item(294, 19).Value = 2.1
$ws.Cells.Item(294, 20).Value = 1.7
$ws.Cells.Item(294, 21).Value = 2.25
$ws.Cells.Item(294, 22).Value = 2.05
$ws.Cells.Item(294, 23).Value = 1.75
$ws.Cells.Item(294, 24).Value = -1
$ws.Cells.Item(294, 25).Value = 1.15
$ws.Cells.Item(294, 27).Value = 0
$ws.Cells.Item(294, 28).Value = 0
$ws.Cells.Item(294, 29).Value = -0.5
$ws.Cells.Item(294, 30).Value = 0.375
